$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 194.81818
$ws.Range("I55").Value = 199.22223
$ws.Range("J55").Value = 175
$ws.Range("K55").Value = 199.22223
$ws.Range("L55").Value = 175
$ws.Range("M55").Value = 14.77777
$ws.Range("N55").Value = -603

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 30627.375
$ws.Range("J31").Value = 111975
$ws.Range("L31").Value = 111975
$ws.Range("N31").Value = -112563

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6416150
$ws.Range("J74").Value = 18395.7
$ws.Range("L74").Value = 18395.7
$ws.Range("N74").Value = -20143.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6416150
$ws.Range("J77").Value = 18395.7
$ws.Range("L77").Value = 91978.5
$ws.Range("N77").Value = -100714.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 9533.825999999999
$ws.Range("J102").Value = 8692.23
$ws.Range("L102").Value = 8692.23
$ws.Range("N102").Value = -11936.23

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 50774.285
$ws.Range("J124").Value = 50774.285
$ws.Range("L124").Value = 50774.285
$ws.Range("N124").Value = -60594.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 48000
$ws.Range("J125").Value = 48000
$ws.Range("L125").Value = 48000
$ws.Range("N125").Value = -57840

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4586.2393
$ws.Range("I132").Value = 3042.0535
$ws.Range("K132").Value = 9126.1605
$ws.Range("M132").Value = -6596.1605

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 58746.75
$ws.Range("J135").Value = 58746.75
$ws.Range("L135").Value = 58746.75
$ws.Range("N135").Value = -68886.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 86999
$ws.Range("J138").Value = 86999
$ws.Range("L138").Value = 86999
$ws.Range("N138").Value = -97279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 95.40000000000001
$ws.Range("I8").Value = 94.888885
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 94.888885
$ws.Range("L8").Value = 100
$ws.Range("M8").Value = 45.111115
$ws.Range("N8").Value = -380

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 39266.09
$ws.Range("J96").Value = 71105.8
$ws.Range("L96").Value = 71105.8
$ws.Range("N96").Value = -76597.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 134000
$ws.Range("J124").Value = 134000
$ws.Range("L124").Value = 134000
$ws.Range("N124").Value = -143820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 59666.668
$ws.Range("J135").Value = 59666.668
$ws.Range("L135").Value = 59666.668
$ws.Range("N135").Value = -69806.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1749.85
$ws.Range("J7").Value = 10947.667
$ws.Range("L7").Value = 10947.667
$ws.Range("N7").Value = -11173.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 532912.25
$ws.Range("I31").Value = 2971.2258
$ws.Range("J31").Value = 1796617.9
$ws.Range("K31").Value = 2971.2258
$ws.Range("L31").Value = 1796617.9
$ws.Range("M31").Value = -2676.2258
$ws.Range("N31").Value = -1797207.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 532912.25
$ws.Range("I34").Value = 2971.2258
$ws.Range("J34").Value = 1796617.9
$ws.Range("K34").Value = 2971.2258
$ws.Range("L34").Value = 1796617.9
$ws.Range("M34").Value = -2769.2258
$ws.Range("N34").Value = -1797021.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 51099.9
$ws.Range("I51").Value = 38500
$ws.Range("J51").Value = 69999.75
$ws.Range("K51").Value = 38500
$ws.Range("L51").Value = 69999.75
$ws.Range("M51").Value = -37764
$ws.Range("N51").Value = -71471.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1614.826
$ws.Range("I58").Value = 1043
$ws.Range("K58").Value = 1043
$ws.Range("M58").Value = -840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 51099.9
$ws.Range("I61").Value = 38500
$ws.Range("J61").Value = 69999.75
$ws.Range("K61").Value = 38500
$ws.Range("L61").Value = 69999.75
$ws.Range("M61").Value = -38152
$ws.Range("N61").Value = -70695.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 44272.43
$ws.Range("J103").Value = 58727.25
$ws.Range("L103").Value = 58727.25
$ws.Range("N103").Value = -61071.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4257.143
$ws.Range("I132").Value = 4257.143
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12771.429
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10241.429
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1614.826
$ws.Range("I136").Value = 1043
$ws.Range("K136").Value = 3129
$ws.Range("M136").Value = -579

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2323.1667
$ws.Range("J68").Value = 2510.75
$ws.Range("L68").Value = 7532.25
$ws.Range("N68").Value = -9154.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2323.1667
$ws.Range("J71").Value = 2510.75
$ws.Range("L71").Value = 22596.75
$ws.Range("N71").Value = -30708.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1487.7646
$ws.Range("I122").Value = 839.6667
$ws.Range("K122").Value = 7557.0003
$ws.Range("M122").Value = -5107.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 100437
$ws.Range("J62").Value = 100437
$ws.Range("L62").Value = 100437
$ws.Range("N62").Value = -101809

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 100437
$ws.Range("J65").Value = 100437
$ws.Range("L65").Value = 301311
$ws.Range("N65").Value = -308175

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 55000
$ws.Range("I88").Value = 55000
$ws.Range("K88").Value = 55000
$ws.Range("M88").Value = -54549

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H91").Value = 55000
$ws.Range("I91").Value = 55000
$ws.Range("K91").Value = 55000
$ws.Range("M91").Value = -53440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3639.2
$ws.Range("I22").Value = 4601
$ws.Range("K22").Value = 4601
$ws.Range("M22").Value = -4306

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3639.2
$ws.Range("I27").Value = 4601
$ws.Range("K27").Value = 4601
$ws.Range("M27").Value = -4494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 47619400
$ws.Range("I55").Value = 62500340
$ws.Range("K55").Value = 62500340
$ws.Range("M55").Value = -62500167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 27399.5
$ws.Range("J95").Value = 27399.5
$ws.Range("L95").Value = 27399.5
$ws.Range("N95").Value = -32891.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6088.7
$ws.Range("I122").Value = 5162.5
$ws.Range("K122").Value = 15487.5
$ws.Range("M122").Value = -13037.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 86780.8
$ws.Range("J127").Value = 86780.8
$ws.Range("L127").Value = 86780.8
$ws.Range("N127").Value = -96700.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 99531.5
$ws.Range("I132").Value = 64105.438
$ws.Range("K132").Value = 192316.314
$ws.Range("M132").Value = -189786.314

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2000684
$ws.Range("I23").Value = 2000684
$ws.Range("K23").Value = 2000684
$ws.Range("M23").Value = -2000455

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10452.566
$ws.Range("I132").Value = 2240.353
$ws.Range("J132").Value = 21191.615
$ws.Range("K132").Value = 6721.059
$ws.Range("L132").Value = 63574.845
$ws.Range("M132").Value = -4191.059
$ws.Range("N132").Value = -68634.845
